$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-key row 3 (Red Pen) onto the existing row-2 style (s=2) by inserting a
# fresh row above it (which inherits row 2's formatting), re-entering the
# values, then dropping the old (now duplicate) row.
$ws.Rows(3).Insert(-4121)
$ws.Range("A3").Value = "Red Pen"
$ws.Range("B3").Value = 200
$ws.Range("C3").Value = 1111144444
$ws.Range("D3").Value = "Office"
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 125
$ws.Rows(4).Delete()

# New header columns G:J
$ws.Range("G1").Value = "Purchase Price "
$ws.Range("H1").Value = "Min Retail Price"
$ws.Range("I1").Value = "Max Retail Price "
$ws.Range("J1").Value = "Min stock Qty"

# New data for row 2 (Blue Pen) - Purchase Price left blank
$ws.Range("H2").Value = 120
$ws.Range("I2").Value = 250
$ws.Range("J2").Value = 10

# New data for row 3 (Red Pen) - Purchase Price left blank
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = 10

# Scroll / selection to match the author's final view
$ws.Range("J3").Select()
